$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) sometimes holds numeric-looking text (e.g. "0.630",
# "35.50", "1.00") that must stay literal text -- otherwise Excel's normal
# "smart" input parsing would coerce it to a Number and silently drop
# meaningful trailing zeros (e.g. "0.630" -> 0.63). Force those specific
# cells to Text format before writing so the digits round-trip exactly.
$textPriceRows = @(4,5,6,7,10,11,12,13,14,17,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,37,38,39,40,43,44,46,48,49,50,51)
foreach ($r in $textPriceRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

# Refresh the cryptocurrency price / 1h-volume data, including the row
# reorder where Bittensor now ranks above OKB.
$ws.Range("D2").Value = '65.437.14'
$ws.Range("E2").Value = '  -0.49%  '
$ws.Range("D3").Value = '3.369.75'
$ws.Range("E3").Value = '  -1.08%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '559.28'
$ws.Range("E5").Value = '  -0.63%  '
$ws.Range("D6").Value = '174.14'
$ws.Range("E6").Value = '  -1.08%  '
$ws.Range("D7").Value = '0.628'
$ws.Range("E7").Value = '  -0.37%  '
$ws.Range("D8").Value = '3.362.45'
$ws.Range("E8").Value = '  -1.01%  '
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("D10").Value = '0.172'
$ws.Range("E10").Value = '  -0.29%  '
$ws.Range("D11").Value = '0.630'
$ws.Range("E11").Value = '  -0.66%  '
$ws.Range("D12").Value = '52.94'
$ws.Range("E12").Value = '  -3.94%  '
$ws.Range("D13").Value = '0.0000276'
$ws.Range("E13").Value = '  -1.87%  '
$ws.Range("D14").Value = '9.15'
$ws.Range("E14").Value = '  -0.37%  '
$ws.Range("D15").Value = '3.892.12'
$ws.Range("E15").Value = '  -1.65%  '
$ws.Range("D17").Value = '18.13'
$ws.Range("E17").Value = '  -1.27%  '
$ws.Range("D18").Value = '3.357.46'
$ws.Range("E18").Value = '  -1.41%  '
$ws.Range("D19").Value = '65.342.31'
$ws.Range("E19").Value = '  -0.40%  '
$ws.Range("D20").Value = '11.77'
$ws.Range("E20").Value = '  -1.04%  '
$ws.Range("D21").Value = '0.991'
$ws.Range("E21").Value = '  -0.23%  '
$ws.Range("D22").Value = '474.23'
$ws.Range("E22").Value = '  +0.56%  '
$ws.Range("D23").Value = '4.90'
$ws.Range("E23").Value = '  -6.75%  '
$ws.Range("D24").Value = '90.06'
$ws.Range("E24").Value = '  +3.95%  '
$ws.Range("D25").Value = '4.07'
$ws.Range("E25").Value = '  -2.15%  '
$ws.Range("D26").Value = '14.13'
$ws.Range("E26").Value = '  +4.36%  '
$ws.Range("D27").Value = '2.88'
$ws.Range("E27").Value = '  -0.49%  '
$ws.Range("D28").Value = '10.52'
$ws.Range("E28").Value = '  -3.54%  '
$ws.Range("D29").Value = '8.64'
$ws.Range("E29").Value = '  -3.06%  '
$ws.Range("D30").Value = '31.05'
$ws.Range("E30").Value = '  -0.02%  '
$ws.Range("D31").Value = '6.49'
$ws.Range("E31").Value = '  -3.24%  '
$ws.Range("D32").Value = '11.37'
$ws.Range("E32").Value = '  -1.86%  '
$ws.Range("B33").Value = 'Bittensor'
$ws.Range("C33").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D33").Value = '574.82'
$ws.Range("E33").Value = '  -0.64%  '
$ws.Range("B34").Value = 'OKB'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D34").Value = '61.87'
$ws.Range("E34").Value = '  -1.40%  '
$ws.Range("E35").Value = '  -1.88%  '
$ws.Range("E36").Value = '  +0.05%  '
$ws.Range("D37").Value = '3.61'
$ws.Range("E37").Value = '  +1.73%  '
$ws.Range("D38").Value = '0.141'
$ws.Range("E38").Value = '  -0.25%  '
$ws.Range("D39").Value = '35.50'
$ws.Range("E39").Value = '  -1.03%  '
$ws.Range("D40").Value = '0.371'
$ws.Range("E40").Value = '  -1.10%  '
$ws.Range("E41").Value = '  -3.66%  '
$ws.Range("D42").Value = '3.087.97'
$ws.Range("E42").Value = '  -0.28%  '
$ws.Range("D43").Value = '2.78'
$ws.Range("E43").Value = '  -2.18%  '
$ws.Range("D44").Value = '0.0413'
$ws.Range("E44").Value = '  -0.76%  '
$ws.Range("E45").Value = '  -2.20%  '
$ws.Range("D46").Value = '3.15'
$ws.Range("E46").Value = '  -1.92%  '
$ws.Range("E47").Value = '  -3.68%  '
$ws.Range("D48").Value = '0.998'
$ws.Range("E48").Value = '  -0.07%  '
$ws.Range("D49").Value = '140.49'
$ws.Range("E49").Value = '  +2.95%  '
$ws.Range("D50").Value = '2.57'
$ws.Range("E50").Value = '  -1.11%  '
$ws.Range("D51").Value = '8.42'
$ws.Range("E51").Value = '  +0.42%  '

Write-Host "Updated cryptos list"
